$d = $word.ActiveDocument

$pairs = @(
    @("957×5=", "302×5="),
    @("752×8=", "823×2="),
    @("559×2=", "663×6="),
    @("114×4=", "187×4="),
    @("845×6=", "377×2="),
    @("582×5=", "235×6="),
    @("172×2=", "331×9="),
    @("967×7=", "515×3="),
    @("195×2=", "553×9="),
    @("367×4=", "657×9="),
    @("287×8=", "829×6="),
    @("194×4=", "520×4="),
    @("939×9=", "504×6="),
    @("363×7=", "516×9="),
    @("324×9=", "780×8="),
    @("136×2=", "276×6="),
    @("800×8=", "635×4="),
    @("944×4=", "424×3="),
    @("629×2=", "337×3="),
    @("199×4=", "758×7="),
    @("259×6=", "401×4="),
    @("506×2=", "704×5="),
    @("335×2=", "480×3="),
    @("575×9=", "116×3="),
    @("653×7=", "603×7=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
